$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NCT(2.3800114516276167, 1.376019635617475, -0.5493296410636643, 2.5855754859949345)"
$ws.Range("C2").Value = "NIG(1.0501168038679465, 0.843168967667512, 3.735646409005841, 5.976638112342744)"
$ws.Range("D2").Value = "NCT(3.0291706242121723, 1.4961017483071146, -0.8874824656087534, 2.3493341870267344)"
$ws.Range("E2").Value = "NIG(1.1960274835144102, 0.7972177607562063, 5.427885264097043, 5.810422201562057)"
